# Updates cryptos.xlsx per commit "Updated symbol list on Thu Jan 19 21:00:06 UTC 2023 with GitHub Actions"
# Columns: D = Price, E = Volume(1h) (percent text), G = Hora (hour, text "20" -> "21")
# All target cells are plain-text strings in the sheet (t="inlineStr"), so each new
# value is written with a leading apostrophe to force Excel to keep it as literal
# text instead of auto-converting it to a Number/Percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Val = '294.35' },
    @{ Addr = "E2"; Val = '1.22%' },
    @{ Addr = "G2"; Val = '21' },
    @{ Addr = "E3"; Val = '0.08%' },
    @{ Addr = "G3"; Val = '21' },
    @{ Addr = "D4"; Val = '4.954' },
    @{ Addr = "E4"; Val = '1.71%' },
    @{ Addr = "G4"; Val = '21' },
    @{ Addr = "D5"; Val = '0.07345' },
    @{ Addr = "E5"; Val = '2.65%' },
    @{ Addr = "G5"; Val = '21' },
    @{ Addr = "D6"; Val = '2.310' },
    @{ Addr = "E6"; Val = '32.66%' },
    @{ Addr = "G6"; Val = '21' },
    @{ Addr = "D7"; Val = '7.742' },
    @{ Addr = "E7"; Val = '0.91%' },
    @{ Addr = "G7"; Val = '21' },
    @{ Addr = "D8"; Val = '3.741' },
    @{ Addr = "E8"; Val = '-0.20%' },
    @{ Addr = "G8"; Val = '21' },
    @{ Addr = "D9"; Val = '0.9079' },
    @{ Addr = "E9"; Val = '1.38%' },
    @{ Addr = "G9"; Val = '21' },
    @{ Addr = "D10"; Val = '0.1681' },
    @{ Addr = "E10"; Val = '2.10%' },
    @{ Addr = "G10"; Val = '21' },
    @{ Addr = "D11"; Val = '0.08072' },
    @{ Addr = "E11"; Val = '9.86%' },
    @{ Addr = "G11"; Val = '21' },
    @{ Addr = "D12"; Val = '0.08171' },
    @{ Addr = "E12"; Val = '2.36%' },
    @{ Addr = "G12"; Val = '21' },
    @{ Addr = "D13"; Val = '0.03103' },
    @{ Addr = "E13"; Val = '3.45%' },
    @{ Addr = "G13"; Val = '21' },
    @{ Addr = "E14"; Val = '0.77%' },
    @{ Addr = "G14"; Val = '21' },
    @{ Addr = "D15"; Val = '0.001508' },
    @{ Addr = "E15"; Val = '1.15%' },
    @{ Addr = "G15"; Val = '21' },
    @{ Addr = "D16"; Val = '0.005755' },
    @{ Addr = "E16"; Val = '1.61%' },
    @{ Addr = "G16"; Val = '21' },
    @{ Addr = "D17"; Val = '3.481' },
    @{ Addr = "E17"; Val = '0.72%' },
    @{ Addr = "G17"; Val = '21' },
    @{ Addr = "D18"; Val = '2.078' },
    @{ Addr = "E18"; Val = '-1.24%' },
    @{ Addr = "G18"; Val = '21' },
    @{ Addr = "E19"; Val = '1.01%' },
    @{ Addr = "G19"; Val = '21' },
    @{ Addr = "E20"; Val = '-0.03%' },
    @{ Addr = "G20"; Val = '21' },
    @{ Addr = "D21"; Val = '3.968' },
    @{ Addr = "E21"; Val = '-9.43%' },
    @{ Addr = "G21"; Val = '21' },
    @{ Addr = "D22"; Val = '0.2098' },
    @{ Addr = "E22"; Val = '4.65%' },
    @{ Addr = "G22"; Val = '21' },
    @{ Addr = "D23"; Val = '0.04550' },
    @{ Addr = "E23"; Val = '1.61%' },
    @{ Addr = "G23"; Val = '21' },
    @{ Addr = "D24"; Val = '0.001212' },
    @{ Addr = "E24"; Val = '-0.30%' },
    @{ Addr = "G24"; Val = '21' },
    @{ Addr = "D25"; Val = '0.004655' },
    @{ Addr = "E25"; Val = '16.46%' },
    @{ Addr = "G25"; Val = '21' },
    @{ Addr = "D26"; Val = '0.0001300' },
    @{ Addr = "E26"; Val = '3.71%' },
    @{ Addr = "G26"; Val = '21' },
    @{ Addr = "D27"; Val = '0.0003395' },
    @{ Addr = "E27"; Val = '-95.49%' },
    @{ Addr = "G27"; Val = '21' },
    @{ Addr = "G28"; Val = '21' },
    @{ Addr = "G29"; Val = '21' },
    @{ Addr = "G30"; Val = '21' },
    @{ Addr = "G31"; Val = '21' },
    @{ Addr = "G32"; Val = '21' },
    @{ Addr = "G33"; Val = '21' },
    @{ Addr = "G34"; Val = '21' },
    @{ Addr = "G35"; Val = '21' },
    @{ Addr = "G36"; Val = '21' },
    @{ Addr = "G37"; Val = '21' },
    @{ Addr = "G38"; Val = '21' },
    @{ Addr = "D39"; Val = '0.01610' },
    @{ Addr = "E39"; Val = '-1.78%' },
    @{ Addr = "G39"; Val = '21' },
    @{ Addr = "D40"; Val = '0.04434' },
    @{ Addr = "E40"; Val = '2.40%' },
    @{ Addr = "G40"; Val = '21' },
    @{ Addr = "D41"; Val = '0.007337' },
    @{ Addr = "E41"; Val = '-0.71%' },
    @{ Addr = "G41"; Val = '21' },
    @{ Addr = "D42"; Val = '0.1332' },
    @{ Addr = "E42"; Val = '1.92%' },
    @{ Addr = "G42"; Val = '21' },
    @{ Addr = "D43"; Val = '0.008605' },
    @{ Addr = "G43"; Val = '21' },
    @{ Addr = "D44"; Val = '0.001947' },
    @{ Addr = "E44"; Val = '-4.71%' },
    @{ Addr = "G44"; Val = '21' },
    @{ Addr = "D45"; Val = '0.009519' },
    @{ Addr = "E45"; Val = '-13.96%' },
    @{ Addr = "G45"; Val = '21' },
    @{ Addr = "D46"; Val = '0.00005729' },
    @{ Addr = "E46"; Val = '-0.16%' },
    @{ Addr = "G46"; Val = '21' },
    @{ Addr = "E47"; Val = '-0.30%' },
    @{ Addr = "G47"; Val = '21' },
    @{ Addr = "E48"; Val = '2.74%' },
    @{ Addr = "G48"; Val = '21' },
    @{ Addr = "D49"; Val = '0.002899' },
    @{ Addr = "E49"; Val = '-3.67%' },
    @{ Addr = "G49"; Val = '21' },
    @{ Addr = "D50"; Val = '0.00002100' },
    @{ Addr = "E50"; Val = '-0.30%' },
    @{ Addr = "G50"; Val = '21' },
    @{ Addr = "D51"; Val = '0.0002000' },
    @{ Addr = "E51"; Val = '-0.30%' },
    @{ Addr = "G51"; Val = '21' }
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = "'" + $u.Val
}
